# Adaptive Right-Sizing deck: enrich "how it works" slides with deeper detail.
$p = $ppt.ActivePresentation

# Slide 1: "Adaptive Warehouse Right-Sizing" overview bullets.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$lines1 = @(
    "Goal: auto-tune warehouse size by observed load patterns",
    "Inputs: staged WAREHOUSE_METERING (credits_used by hour)",
    "Policy DT: RIGHT_SIZING_POLICY_DT → per-warehouse/hour recommendation",
    "Executor: APPLY_RIGHT_SIZING() applies size + optional multi-cluster"
)
$tr1.Text = [string]::Join("`r", $lines1)

# Slide 2: "How it works" deep-dive bullets (grows from 4 to 6 bullets).
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$lines2 = @(
    "Ingestion: Task merges ACCOUNT_USAGE.WAREHOUSE_METERING → TECHUP.AUDIT.WAREHOUSE_METERING_STG (change tracking)",
    "Signal: Aggregate credits_used into hourly buckets per warehouse",
    "Policy logic: map avg(credits_used) ranges → SMALL/MEDIUM/LARGE sizing and multi-cluster toggle",
    "Governance: all changes logged in RIGHT_SIZING_LOG with status, DDL, error",
    "Orchestration: APPLY_RIGHT_SIZING_TASK executes on-the-hour against current hour recommendation",
    "Safety: thresholds are conservative; tune sizing cutoffs per environment; dry-run by commenting execute immediate"
)
$tr2.Text = [string]::Join("`r", $lines2)
